$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column E ("purpose") rows 2-19 contain "fullRNASEQ" and should read "fullRNASeq".
for ($row = 2; $row -le 19; $row++) {
    $ws.Cells.Item($row, 5).Value = "fullRNASeq"  # Column E
}
